# Update computed market/profit columns (H-N) on several Leve tables.
# Values below were refreshed from the latest market-board snapshot;
# row/column coordinates and the set of touched cells are unchanged.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 32
$ws.Cells.Item(32, 8).Value = 325.5  # H32: 655.1667 -> 325.5
$ws.Cells.Item(32, 9).Value = 100  # I32: 600 -> 100
$ws.Cells.Item(32, 10).Value = 1002  # J32: 673.55554 -> 1002
$ws.Cells.Item(32, 11).Value = 100  # K32: 600 -> 100
$ws.Cells.Item(32, 12).Value = 1002  # L32: 673.55554 -> 1002
$ws.Cells.Item(32, 13).Value = 226  # M32: -274 -> 226
$ws.Cells.Item(32, 14).Value = -1654  # N32: -1325.55554 -> -1654

# ALC!row 129
$ws.Cells.Item(129, 8).Value = 2552105.5  # H129: 3206270 -> 2552105.5
$ws.Cells.Item(129, 9).Value = 27778752  # I129: 31251072 -> 27778752
$ws.Cells.Item(129, 10).Value = 1096.3707  # J129: 1149.9572 -> 1096.3707
$ws.Cells.Item(129, 11).Value = 83336256  # K129: 93753216 -> 83336256
$ws.Cells.Item(129, 12).Value = 3289.1121  # L129: 3449.8716 -> 3289.1121
$ws.Cells.Item(129, 13).Value = -83331256  # M129: -93748216 -> -83331256
$ws.Cells.Item(129, 14).Value = -13289.1121  # N129: -13449.8716 -> -13289.1121

# ALC!row 137
$ws.Cells.Item(137, 8).Value = 2566817.5  # H137: 3229004.5 -> 2566817.5
$ws.Cells.Item(137, 9).Value = 3450996.5  # I137: 4765340.5 -> 3450996.5
$ws.Cells.Item(137, 11).Value = 10352989.5  # K137: 14296021.5 -> 10352989.5
$ws.Cells.Item(137, 13).Value = -10350439.5  # M137: -14293471.5 -> -10350439.5

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 32
$ws.Cells.Item(32, 8).Value = 3072.47  # H32: 3657.75 -> 3072.47
$ws.Cells.Item(32, 9).Value = 2640.6353  # I32: 3060.7805 -> 2640.6353
$ws.Cells.Item(32, 10).Value = 5519.533  # J32: 6377.278 -> 5519.533
$ws.Cells.Item(32, 11).Value = 2640.6353  # K32: 3060.7805 -> 2640.6353
$ws.Cells.Item(32, 12).Value = 5519.533  # L32: 6377.278 -> 5519.533
$ws.Cells.Item(32, 13).Value = -2353.6353  # M32: -2773.7805 -> -2353.6353
$ws.Cells.Item(32, 14).Value = -6093.533  # N32: -6951.278 -> -6093.533

# ARM!row 61
$ws.Cells.Item(61, 8).Value = 1228.9348  # H61: 1476.94 -> 1228.9348
$ws.Cells.Item(61, 9).Value = 631.6905  # I61: 639.4878 -> 631.6905
$ws.Cells.Item(61, 10).Value = 7500  # J61: 5292 -> 7500
$ws.Cells.Item(61, 11).Value = 631.6905  # K61: 639.4878 -> 631.6905
$ws.Cells.Item(61, 12).Value = 7500  # L61: 5292 -> 7500
$ws.Cells.Item(61, 13).Value = -419.6905  # M61: -427.4878 -> -419.6905
$ws.Cells.Item(61, 14).Value = -7924  # N61: -5716 -> -7924

# ARM!row 74
$ws.Cells.Item(74, 8).Value = 685.55554  # H74: 701.64703 -> 685.55554
$ws.Cells.Item(74, 9).Value = 685.55554  # I74: 701.64703 -> 685.55554
$ws.Cells.Item(74, 11).Value = 685.55554  # K74: 701.64703 -> 685.55554
$ws.Cells.Item(74, 13).Value = 188.44446  # M74: 172.35297 -> 188.44446

# ARM!row 77
$ws.Cells.Item(77, 8).Value = 685.55554  # H77: 701.64703 -> 685.55554
$ws.Cells.Item(77, 9).Value = 685.55554  # I77: 701.64703 -> 685.55554
$ws.Cells.Item(77, 11).Value = 3427.7777  # K77: 3508.23515 -> 3427.7777
$ws.Cells.Item(77, 13).Value = 940.2223000000004  # M77: 859.76485 -> 940.2223000000004

# ARM!row 132
$ws.Cells.Item(132, 8).Value = 1985.2115  # H132: 2065.1 -> 1985.2115
$ws.Cells.Item(132, 9).Value = 1574.1316  # I132: 1549.7693 -> 1574.1316
$ws.Cells.Item(132, 10).Value = 3101  # J132: 3892.182 -> 3101
$ws.Cells.Item(132, 11).Value = 4722.3948  # K132: 4649.3079 -> 4722.3948
$ws.Cells.Item(132, 12).Value = 9303  # L132: 11676.546 -> 9303
$ws.Cells.Item(132, 13).Value = -2192.3948  # M132: -2119.3079 -> -2192.3948
$ws.Cells.Item(132, 14).Value = -14363  # N132: -16736.546 -> -14363

# ARM!row 133
$ws.Cells.Item(133, 8).Value = 30000  # H133: 26653.666 -> 30000
$ws.Cells.Item(133, 10).Value = 30000  # J133: 26653.666 -> 30000
$ws.Cells.Item(133, 12).Value = 30000  # L133: 26653.666 -> 30000
$ws.Cells.Item(133, 14).Value = -35060  # N133: -31713.666 -> -35060

# ARM!row 136
$ws.Cells.Item(136, 8).Value = 1228.9348  # H136: 1476.94 -> 1228.9348
$ws.Cells.Item(136, 9).Value = 631.6905  # I136: 639.4878 -> 631.6905
$ws.Cells.Item(136, 10).Value = 7500  # J136: 5292 -> 7500
$ws.Cells.Item(136, 11).Value = 1895.0715  # K136: 1918.4634 -> 1895.0715
$ws.Cells.Item(136, 12).Value = 22500  # L136: 15876 -> 22500
$ws.Cells.Item(136, 13).Value = 654.9285  # M136: 631.5365999999999 -> 654.9285
$ws.Cells.Item(136, 14).Value = -27600  # N136: -20976 -> -27600

$ws = $wb.Worksheets.Item("BSM")
# BSM!row 105
$ws.Cells.Item(105, 8).Value = 1524.3684  # H105: 1597.3334 -> 1524.3684
$ws.Cells.Item(105, 9).Value = 1359.2858  # I105: 1430.7142 -> 1359.2858
$ws.Cells.Item(105, 10).Value = 1986.6  # J105: 2180.5 -> 1986.6
$ws.Cells.Item(105, 11).Value = 1359.2858  # K105: 1430.7142 -> 1359.2858
$ws.Cells.Item(105, 12).Value = 1986.6  # L105: 2180.5 -> 1986.6
$ws.Cells.Item(105, 13).Value = 387.7141999999999  # M105: 316.2858000000001 -> 387.7141999999999
$ws.Cells.Item(105, 14).Value = -5480.6  # N105: -5674.5 -> -5480.6

# BSM!row 134
$ws.Cells.Item(134, 8).Value = 1918.2759  # H134: 1990.5264 -> 1918.2759
$ws.Cells.Item(134, 9).Value = 1446.6666  # I134: 1524.2554 -> 1446.6666
$ws.Cells.Item(134, 11).Value = 4339.9998  # K134: 4572.7662 -> 4339.9998
$ws.Cells.Item(134, 13).Value = -1804.9998  # M134: -2037.7662 -> -1804.9998

$ws = $wb.Worksheets.Item("CRP")
# CRP!row 31
$ws.Cells.Item(31, 8).Value = 2383798.5  # H31: 2860423.5 -> 2383798.5
$ws.Cells.Item(31, 9).Value = 3847664.2  # I31: 4763617 -> 3847664.2
$ws.Cells.Item(31, 10).Value = 5016.3125  # J31: 5632.9287 -> 5016.3125
$ws.Cells.Item(31, 11).Value = 3847664.2  # K31: 4763617 -> 3847664.2
$ws.Cells.Item(31, 12).Value = 5016.3125  # L31: 5632.9287 -> 5016.3125
$ws.Cells.Item(31, 13).Value = -3847369.2  # M31: -4763322 -> -3847369.2
$ws.Cells.Item(31, 14).Value = -5606.3125  # N31: -6222.9287 -> -5606.3125

# CRP!row 34
$ws.Cells.Item(34, 8).Value = 2383798.5  # H34: 2860423.5 -> 2383798.5
$ws.Cells.Item(34, 9).Value = 3847664.2  # I34: 4763617 -> 3847664.2
$ws.Cells.Item(34, 10).Value = 5016.3125  # J34: 5632.9287 -> 5016.3125
$ws.Cells.Item(34, 11).Value = 3847664.2  # K34: 4763617 -> 3847664.2
$ws.Cells.Item(34, 12).Value = 5016.3125  # L34: 5632.9287 -> 5016.3125
$ws.Cells.Item(34, 13).Value = -3847462.2  # M34: -4763415 -> -3847462.2
$ws.Cells.Item(34, 14).Value = -5420.3125  # N34: -6036.9287 -> -5420.3125

# CRP!row 58
$ws.Cells.Item(58, 8).Value = 8623151  # H58: 9261869 -> 8623151
$ws.Cells.Item(58, 9).Value = 1533.4524  # I58: 1590.125 -> 1533.4524
$ws.Cells.Item(58, 10).Value = 31254896  # J58: 35719810 -> 31254896
$ws.Cells.Item(58, 11).Value = 1533.4524  # K58: 1590.125 -> 1533.4524
$ws.Cells.Item(58, 12).Value = 31254896  # L58: 35719810 -> 31254896
$ws.Cells.Item(58, 13).Value = -1330.4524  # M58: -1387.125 -> -1330.4524
$ws.Cells.Item(58, 14).Value = -31255302  # N58: -35720216 -> -31255302

# CRP!row 132
$ws.Cells.Item(132, 8).Value = 1608.551  # H132: 1668.4894 -> 1608.551
$ws.Cells.Item(132, 9).Value = 1196.1464  # I132: 1223.55 -> 1196.1464
$ws.Cells.Item(132, 10).Value = 3722.125  # J132: 4211 -> 3722.125
$ws.Cells.Item(132, 11).Value = 3588.4392  # K132: 3670.65 -> 3588.4392
$ws.Cells.Item(132, 12).Value = 11166.375  # L132: 12633 -> 11166.375
$ws.Cells.Item(132, 13).Value = -1058.4392  # M132: -1140.65 -> -1058.4392
$ws.Cells.Item(132, 14).Value = -16226.375  # N132: -17693 -> -16226.375

# CRP!row 134
$ws.Cells.Item(134, 8).Value = 1563.875  # H134: 1503.262 -> 1563.875
$ws.Cells.Item(134, 9).Value = 862.21875  # I134: 828.6177 -> 862.21875
$ws.Cells.Item(134, 11).Value = 2586.65625  # K134: 2485.8531 -> 2586.65625
$ws.Cells.Item(134, 13).Value = -51.65625  # M134: 49.14689999999973 -> -51.65625

# CRP!row 136
$ws.Cells.Item(136, 8).Value = 8623151  # H136: 9261869 -> 8623151
$ws.Cells.Item(136, 9).Value = 1533.4524  # I136: 1590.125 -> 1533.4524
$ws.Cells.Item(136, 10).Value = 31254896  # J136: 35719810 -> 31254896
$ws.Cells.Item(136, 11).Value = 4600.357199999999  # K136: 4770.375 -> 4600.357199999999
$ws.Cells.Item(136, 12).Value = 93764688  # L136: 107159430 -> 93764688
$ws.Cells.Item(136, 13).Value = -2050.357199999999  # M136: -2220.375 -> -2050.357199999999
$ws.Cells.Item(136, 14).Value = -93769788  # N136: -107164530 -> -93769788

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 64
$ws.Cells.Item(64, 8).Value = 14262.625  # H64: 25650.25 -> 14262.625
$ws.Cells.Item(64, 10).Value = 22300  # J64: 100000 -> 22300
$ws.Cells.Item(64, 12).Value = 66900  # L64: 300000 -> 66900
$ws.Cells.Item(64, 14).Value = -67440  # N64: -300540 -> -67440

# CUL!row 67
$ws.Cells.Item(67, 8).Value = 14262.625  # H67: 25650.25 -> 14262.625
$ws.Cells.Item(67, 10).Value = 22300  # J67: 100000 -> 22300
$ws.Cells.Item(67, 12).Value = 66900  # L67: 300000 -> 66900
$ws.Cells.Item(67, 14).Value = -68772  # N67: -301872 -> -68772

# CUL!row 87
$ws.Cells.Item(87, 8).Value = 11207.895  # H87: 11692.777 -> 11207.895
$ws.Cells.Item(87, 9).Value = 9925  # I87: 10497.692 -> 9925
$ws.Cells.Item(87, 11).Value = 29775  # K87: 31493.076 -> 29775
$ws.Cells.Item(87, 13).Value = -28527  # M87: -30245.076 -> -28527

# CUL!row 90
$ws.Cells.Item(90, 8).Value = 11207.895  # H90: 11692.777 -> 11207.895
$ws.Cells.Item(90, 9).Value = 9925  # I90: 10497.692 -> 9925
$ws.Cells.Item(90, 11).Value = 89325  # K90: 94479.22799999999 -> 89325
$ws.Cells.Item(90, 13).Value = -83085  # M90: -88239.22799999999 -> -83085

# CUL!row 107
$ws.Cells.Item(107, 8).Value = 699.90247  # H107: 710.0244 -> 699.90247
$ws.Cells.Item(107, 9).Value = 237.77777  # I107: 258.3 -> 237.77777
$ws.Cells.Item(107, 10).Value = 829.875  # J107: 855.74194 -> 829.875
$ws.Cells.Item(107, 11).Value = 713.33331  # K107: 774.9000000000001 -> 713.33331
$ws.Cells.Item(107, 12).Value = 2489.625  # L107: 2567.22582 -> 2489.625
$ws.Cells.Item(107, 13).Value = 1206.66669  # M107: 1145.1 -> 1206.66669
$ws.Cells.Item(107, 14).Value = -6329.625  # N107: -6407.22582 -> -6329.625

# CUL!row 133
$ws.Cells.Item(133, 8).Value = 4177.6113  # H133: 3840.2856 -> 4177.6113
$ws.Cells.Item(133, 9).Value = 5808.5713  # I133: 4761.1113 -> 5808.5713
$ws.Cells.Item(133, 10).Value = 3139.7273  # J133: 3149.6667 -> 3139.7273
$ws.Cells.Item(133, 11).Value = 17425.7139  # K133: 14283.3339 -> 17425.7139
$ws.Cells.Item(133, 12).Value = 9419.1819  # L133: 9449.000100000001 -> 9419.1819
$ws.Cells.Item(133, 13).Value = -12365.7139  # M133: -9223.333899999998 -> -12365.7139
$ws.Cells.Item(133, 14).Value = -19539.1819  # N133: -19569.0001 -> -19539.1819

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 45
$ws.Cells.Item(45, 8).Value = 6975  # H45: 5833.1665 -> 6975
$ws.Cells.Item(45, 9).Value = 4950  # I45: 4249.75 -> 4950
$ws.Cells.Item(45, 11).Value = 4950  # K45: 4249.75 -> 4950
$ws.Cells.Item(45, 13).Value = -4543  # M45: -3842.75 -> -4543

# LTW!row 46
$ws.Cells.Item(46, 8).Value = 2330  # H46: 2563.7273 -> 2330
$ws.Cells.Item(46, 9).Value = 600  # I46: 533.6667 -> 600
$ws.Cells.Item(46, 10).Value = 2676  # J46: 3325 -> 2676
$ws.Cells.Item(46, 11).Value = 600  # K46: 533.6667 -> 600
$ws.Cells.Item(46, 12).Value = 2676  # L46: 3325 -> 2676
$ws.Cells.Item(46, 13).Value = -412  # M46: -345.6667 -> -412
$ws.Cells.Item(46, 14).Value = -3052  # N46: -3701 -> -3052

# LTW!row 68
$ws.Cells.Item(68, 8).Value = 2171.0967  # H68: 2210.1667 -> 2171.0967
$ws.Cells.Item(68, 9).Value = 1060.04  # I68: 1062.5834 -> 1060.04
$ws.Cells.Item(68, 11).Value = 1060.04  # K68: 1062.5834 -> 1060.04
$ws.Cells.Item(68, 13).Value = -311.04  # M68: -313.5834 -> -311.04

# LTW!row 71
$ws.Cells.Item(71, 8).Value = 2171.0967  # H71: 2210.1667 -> 2171.0967
$ws.Cells.Item(71, 9).Value = 1060.04  # I71: 1062.5834 -> 1060.04
$ws.Cells.Item(71, 11).Value = 5300.2  # K71: 5312.916999999999 -> 5300.2
$ws.Cells.Item(71, 13).Value = -1556.2  # M71: -1568.916999999999 -> -1556.2

# LTW!row 132
$ws.Cells.Item(132, 8).Value = 1839.6888  # H132: 2132.2703 -> 1839.6888
$ws.Cells.Item(132, 9).Value = 1231.8064  # I132: 1411.76 -> 1231.8064
$ws.Cells.Item(132, 10).Value = 3185.7144  # J132: 3633.3333 -> 3185.7144
$ws.Cells.Item(132, 11).Value = 3695.4192  # K132: 4235.28 -> 3695.4192
$ws.Cells.Item(132, 12).Value = 9557.143199999999  # L132: 10899.9999 -> 9557.143199999999
$ws.Cells.Item(132, 13).Value = -1165.4192  # M132: -1705.28 -> -1165.4192
$ws.Cells.Item(132, 14).Value = -14617.1432  # N132: -15959.9999 -> -14617.1432

# LTW!row 133
$ws.Cells.Item(133, 8).Value = 27963  # H133: 42000 -> 27963
$ws.Cells.Item(133, 10).Value = 27963  # J133: 42000 -> 27963
$ws.Cells.Item(133, 12).Value = 27963  # L133: 42000 -> 27963
$ws.Cells.Item(133, 14).Value = -33023  # N133: -47060 -> -33023

# LTW!row 136
$ws.Cells.Item(136, 8).Value = 2566099  # H136: 2440944.2 -> 2566099
$ws.Cells.Item(136, 9).Value = 3126511.2  # I136: 2942624 -> 3126511.2
$ws.Cells.Item(136, 11).Value = 9379533.600000001  # K136: 8827872 -> 9379533.600000001
$ws.Cells.Item(136, 13).Value = -9376983.600000001  # M136: -8825322 -> -9376983.600000001

$ws = $wb.Worksheets.Item("WVR")
# WVR!row 62
$ws.Cells.Item(62, 8).Value = 3800  # H62: 3117.2144 -> 3800
$ws.Cells.Item(62, 9).Value = 3800  # I62: 3015.4614 -> 3800
$ws.Cells.Item(62, 10).Value = 0  # J62: 3205.4 -> 0
$ws.Cells.Item(62, 11).Value = 3800  # K62: 3015.4614 -> 3800
$ws.Cells.Item(62, 12).Value = 0  # L62: 3205.4 -> 0
$ws.Cells.Item(62, 13).Value = -3176  # M62: -2391.4614 -> -3176
$ws.Cells.Item(62, 14).Value = $null  # N62: cleared

# WVR!row 65
$ws.Cells.Item(65, 8).Value = 3800  # H65: 3117.2144 -> 3800
$ws.Cells.Item(65, 9).Value = 3800  # I65: 3015.4614 -> 3800
$ws.Cells.Item(65, 10).Value = 0  # J65: 3205.4 -> 0
$ws.Cells.Item(65, 11).Value = 19000  # K65: 15077.307 -> 19000
$ws.Cells.Item(65, 12).Value = 0  # L65: 16027 -> 0
$ws.Cells.Item(65, 13).Value = -15880  # M65: -11957.307 -> -15880
$ws.Cells.Item(65, 14).Value = $null  # N65: cleared

# WVR!row 132
$ws.Cells.Item(132, 8).Value = 189630.47  # H132: 216691.88 -> 189630.47
$ws.Cells.Item(132, 9).Value = 234635.05  # I132: 265497.44 -> 234635.05
$ws.Cells.Item(132, 10).Value = 40769.152  # J132: 48090.816 -> 40769.152
$ws.Cells.Item(132, 11).Value = 703905.1499999999  # K132: 796492.3200000001 -> 703905.1499999999
$ws.Cells.Item(132, 12).Value = 122307.456  # L132: 144272.448 -> 122307.456
$ws.Cells.Item(132, 13).Value = -701375.1499999999  # M132: -793962.3200000001 -> -701375.1499999999
$ws.Cells.Item(132, 14).Value = -127367.456  # N132: -149332.448 -> -127367.456

# WVR!row 136
$ws.Cells.Item(136, 8).Value = 1655.381  # H136: 1521.68 -> 1655.381
$ws.Cells.Item(136, 9).Value = 1151.3334  # I136: 1017 -> 1151.3334
$ws.Cells.Item(136, 10).Value = 1857  # J136: 1918.2142 -> 1857
$ws.Cells.Item(136, 11).Value = 3454.0002  # K136: 3051 -> 3454.0002
$ws.Cells.Item(136, 12).Value = 5571  # L136: 5754.642599999999 -> 5571
$ws.Cells.Item(136, 13).Value = -904.0001999999999  # M136: -501 -> -904.0001999999999
$ws.Cells.Item(136, 14).Value = -10671  # N136: -10854.6426 -> -10671
